# Update countries & provincias Spain
# --------------------------------------------------------------
# This script reproduces a data refresh of the COVID "Pais" sheet:
#   1. Two countries (Maldivas, Uganda) moved up one rank in the
#      table (their totals overtook the country previously just
#      above them), which cascades the totals of the countries
#      between the old and new rank down by one row, and gives
#      the promoted country brand-new totals.
#   2. A handful of unrelated rows got refreshed totals (USA,
#      Arabia Saudita, San Marino, Liechtenstein block, etc.)
#   3. The "last updated" timestamp advanced from 23:22 to 23:52.
# --------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Row 4 (Estados Unidos) ----
$ws.Cells.Item(4, 2).Value = 1092492
$ws.Cells.Item(4, 3).Value = 28298
$ws.Cells.Item(4, 4).Value = 151774
$ws.Cells.Item(4, 5).Value = 876955
$ws.Cells.Item(4, 7).Value = 2108
$ws.Cells.Item(4, 8).Value = 63763

# ---- 2. Row 23 (Arabia Saudita) ----
$ws.Cells.Item(23, 6).Value = 123

# ---- 3. Row 104 (San Marino) ----
$ws.Cells.Item(104, 2).Value = 645
$ws.Cells.Item(104, 3).Value = 4
$ws.Cells.Item(104, 4).Value = 506
$ws.Cells.Item(104, 5).Value = 96

# ---- 4. Maldivas overtakes Malta: row 114 gets Maldivas' new
#         totals, rows 115-119 cascade down by one (each row takes
#         on the values the row above used to hold). ----
$ws.Cells.Item(114, 1).Value = "Maldivas"
$ws.Cells.Item(115, 1).Value = "Malta"
$ws.Cells.Item(116, 1).Value = "Jordania"
$ws.Cells.Item(117, 1).Value = "Sudan"
$ws.Cells.Item(118, 1).Value = "Taiwan"
$ws.Cells.Item(119, 1).Value = "Reunion"

$ws.Cells.Item(114, 2).Value = 468
$ws.Cells.Item(114, 3).Value = 190
$ws.Cells.Item(114, 4).Value = 17
$ws.Cells.Item(114, 5).Value = 450
$ws.Cells.Item(114, 6).Value = 2
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 1

$ws.Cells.Item(115, 2).Value = 465
$ws.Cells.Item(115, 3).Value = 2
$ws.Cells.Item(115, 4).Value = 351
$ws.Cells.Item(115, 5).Value = 110
$ws.Cells.Item(115, 6).Value = 1
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 4

$ws.Cells.Item(116, 2).Value = 453
$ws.Cells.Item(116, 3).Value = 2
$ws.Cells.Item(116, 4).Value = 362
$ws.Cells.Item(116, 5).Value = 83
$ws.Cells.Item(116, 6).Value = 5
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 8

$ws.Cells.Item(117, 2).Value = 442
$ws.Cells.Item(117, 3).Value = 67
$ws.Cells.Item(117, 4).Value = 39
$ws.Cells.Item(117, 5).Value = 372
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 3
$ws.Cells.Item(117, 8).Value = 31

$ws.Cells.Item(118, 2).Value = 429
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 4).Value = 322
$ws.Cells.Item(118, 5).Value = 101
$ws.Cells.Item(118, 6).Value = 0
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 6

$ws.Cells.Item(119, 2).Value = 420
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 300
$ws.Cells.Item(119, 5).Value = 120
$ws.Cells.Item(119, 6).Value = 2
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 0

# ---- 5. Row 123 (Estado de Palestina) ----
$ws.Cells.Item(123, 4).Value = 76
$ws.Cells.Item(123, 5).Value = 266

# ---- 6. Uganda overtakes Liechtenstein: row 155 gets Uganda's new
#         totals, rows 156-157 cascade down by one. ----
$ws.Cells.Item(155, 1).Value = "Uganda"
$ws.Cells.Item(156, 1).Value = "Liechtenstein"
$ws.Cells.Item(157, 1).Value = "Barbados"

$ws.Cells.Item(155, 2).Value = 83
$ws.Cells.Item(155, 3).Value = 2
$ws.Cells.Item(155, 4).Value = 52
$ws.Cells.Item(155, 5).Value = 31
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 0

$ws.Cells.Item(156, 2).Value = 82
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 55
$ws.Cells.Item(156, 5).Value = 26
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 1

$ws.Cells.Item(157, 3).Value = 1
$ws.Cells.Item(157, 4).Value = 39
$ws.Cells.Item(157, 5).Value = 35
$ws.Cells.Item(157, 6).Value = 4
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 7

# ---- 7. Update the "last updated" timestamp (row 1) ----
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 23:52"
